$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPSC")

# Add the new "importance" header in K1, matching the style used by the
# other header cells (J1, etc.)
$ws.Range("K1").Value = "importance"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Importance values for courses in rows 2-71 (row 72 "CPSC 3XXX" has none)
$importance = @(
    100, 95, 90, -1, 85, 65, 20, 20, 50, 25,
    20, -1, 10, 50, 45, 50, -1, 20, 40, 75,
    20, 35, 40, 30, 25, -1, 25, 30, 25, 35,
    25, 25, 15, 25, 25, 25, 25, 25, 25, 35,
    30, 25, 25, 25, 25, 25, 25, 65, -1, 20,
    15, 20, 5, 35, 35, 15, 15, -1, 15, 15,
    -1, 15, 20, -1, 20, 20, 25, -1, -1, -1
)

for ($i = 0; $i -lt $importance.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $importance[$i]
}
